# Update the division-problem worksheet table.
# The table has 20 rows x 5 columns; only rows 1, 5, 9, 13, 17 (1-indexed)
# contain problems, the rest are blank answer rows.
# Cell.Range.Text replaces the run text directly.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

function Set-Cell($row, $col, $value) {
    $cell = $t.Cell($row, $col)
    $cell.Range.Text = $value
}

Set-Cell 1 1 "53÷9="
Set-Cell 1 2 "54÷7="
Set-Cell 1 3 "62÷4="
Set-Cell 1 4 "63÷6="
Set-Cell 1 5 "17÷6="

Set-Cell 5 1 "12÷8="
Set-Cell 5 2 "41÷2="
Set-Cell 5 3 "29÷8="
Set-Cell 5 4 "64÷8="
Set-Cell 5 5 "96÷3="

Set-Cell 9 1 "72÷7="
Set-Cell 9 2 "90÷7="
Set-Cell 9 3 "69÷4="
Set-Cell 9 4 "91÷4="
Set-Cell 9 5 "20÷5="

Set-Cell 13 1 "22÷2="
# Row 13, Col 2 ("34÷8=") is unchanged per the diff.
Set-Cell 13 3 "35÷2="
Set-Cell 13 4 "27÷2="
Set-Cell 13 5 "18÷6="

Set-Cell 17 1 "54÷3="
Set-Cell 17 2 "64÷4="
Set-Cell 17 3 "58÷2="
Set-Cell 17 4 "25÷7="
Set-Cell 17 5 "20÷8="

Write-Host "Done"
